$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting existing rows 33-36 down to 34-37
$ws.Rows.Item(33).Insert()

# Fill in the new row 33 with the latest weekly price record
$ws.Cells.Item(33,1).Value = 1
$ws.Cells.Item(33,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(33,3).Value = "Arica y Parinacota"
$ws.Cells.Item(33,4).Value = "2023-04-05"
$ws.Cells.Item(33,5).Value = 15
$ws.Cells.Item(33,6).Value = "Fruta"
$ws.Cells.Item(33,7).Value = 100103
$ws.Cells.Item(33,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(33,9).Value = 100103002
$ws.Cells.Item(33,10).Value = "Ciruela"
$ws.Cells.Item(33,11).Value = "Angeleno"
$ws.Cells.Item(33,12).Value = "Segunda"
$ws.Cells.Item(33,13).Value = 300
$ws.Cells.Item(33,14).Value = 19000
$ws.Cells.Item(33,15).Value = 20000
$ws.Cells.Item(33,16).Value = 19500
$ws.Cells.Item(33,17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(33,18).Value = "Región de O'Higgins"
$ws.Cells.Item(33,19).Value = 1083
$ws.Cells.Item(33,20).Value = 18
